$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original Nocturne de Lure entry (row 13) is now marked as postponed to 11 September
$ws.Range("G13").Value = "Reportée au 11 septembre"

# Insert a new row at position 46 (shifts existing rows 46+ down by one)
$ws.Rows("46:46").Insert()

# Fill in the new row 46 with the "new date" entry for the Nocturne de Lure event
$ws.Range("A46").Value = "Ven 11 Septembre"
$ws.Range("B46").Value = "Nocturne de Lure (ouvert aux cadets)  "
$ws.Range("C46").Value = "VC Luron"
$ws.Range("D46").Value = "Route"
$ws.Range("G46").Value = "Nouvelle date"

# Row 47 (formerly row 46, "Les 3h VTT du VCSA") also gets the "Nouvelle date" marker
$ws.Range("G47").Value = "Nouvelle date"

$ws.Range("E46").Value = "lure_old"

# Update view selection to reflect the author's final cursor position
$ws.Activate()
$ws.Range("E47").Select()
